# actional.xlsx: "jsbr ET 4h 33m"
# The execution-time entry for the "js in browser" action (cell C2) is
# corrected from "4h 3m" to "4h 33m".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "4h 33m"

# Leave the cursor/selection on the cell that was just edited, matching
# the saved view state of the edited workbook.
$ws.Range("C2").Select() | Out-Null
